# Insert 6 new data rows before the current row 267, shifting all rows from
# 267 downward to 273 onward. This mirrors the OOXML diff, where rows
# 267-350 end up at 273-357 and 6 brand-new rows (with a new reporting date,
# 2021-11-10 / serial 44510) are introduced at 267-272.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A267:T272").Insert()

# Constant values shared by every data row in this sheet.
$mercadoId = 2
$mercado   = "Comercializadora del Agro de Limarí"
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100102
$producto   = "Cítricos"
$categoriaId = 100102005
$categoria   = "Naranja"
$unidad    = "`$/bins (400 kilos)"
$origen    = "Provincia de Limarí"
$kgUnidad  = 400

$fecha = 44510

# K, L, M, N, O, P, S for each of the 6 new rows (267-272)
$newRows = @(
    @{ Row = 267; Variedad = "Cara cara";  Calidad = "Primera"; Volumen = 20; PMin = 175000; PMax = 180000; PProm = 177500; PKg = 444 },
    @{ Row = 268; Variedad = "Cara cara";  Calidad = "Segunda"; Volumen = 20; PMin = 155000; PMax = 160000; PProm = 157500; PKg = 394 },
    @{ Row = 269; Variedad = "Lane Late";  Calidad = "Primera"; Volumen = 20; PMin = 155000; PMax = 160000; PProm = 157500; PKg = 394 },
    @{ Row = 270; Variedad = "Lane Late";  Calidad = "Segunda"; Volumen = 20; PMin = 125000; PMax = 130000; PProm = 127500; PKg = 319 },
    @{ Row = 271; Variedad = "Navel Late"; Calidad = "Primera"; Volumen = 20; PMin = 155000; PMax = 160000; PProm = 157500; PKg = 394 },
    @{ Row = 272; Variedad = "Navel Late"; Calidad = "Segunda"; Volumen = 16; PMin = 125000; PMax = 130000; PProm = 127500; PKg = 319 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $r.Variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}
